$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (hyperlink cell style) used by existing link cells in column A
$linkStyle = $ws.Range("A53").Style

# --- New row 54: genomeweb.com copy of the Biocartis article ---
$url54 = "https://www.genomeweb.com/companion-diagnostics/biocartis-positioning-idylla-system-sample-answer-oncology-cdx-tests"
$title = "Biocartis Positioning Idylla System for Sample-to-Answer Oncology CDx Tests"

$ws.Range("A54").Value = $url54
$ws.Range("B54").Value = "CDx"
$ws.Range("C54").Value = $title
$ws.Hyperlinks.Add($ws.Range("A54"), $url54)
$ws.Range("A54").Style = $linkStyle

# --- New row 55: 360dx.com copy of the same Biocartis article ---
$url55 = "https://www.360dx.com/companion-diagnostics/biocartis-positioning-idylla-system-sample-answer-oncology-cdx-tests"

$ws.Range("A55").Value = $url55
$ws.Range("B55").Value = "CDx"
$ws.Range("C55").Value = $title
$ws.Hyperlinks.Add($ws.Range("A55"), $url55)
$ws.Range("A55").Style = $linkStyle
